$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("M2").Value = 24.75542533333333
$ws.Range("N2").Value = 74.26627599999999
$ws.Range("O2").Value = 0.7762421087066456
$ws.Range("P2").Value = 0.7762421087066456
$ws.Range("Q2").Value = 13.32943499360666
$ws.Range("R2").Value = 119.96491494246
$ws.Range("S2").Value = 0.02617184105339938
$ws.Range("T2").Value = 0.02617184105339938
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("M3").Value = 3.818542
$ws.Range("O3").Value = 0.1197358984688377
$ws.Range("P3").Value = 0.1197358984688377
$ws.Range("Q3").Value = 2.05607484719
$ws.Range("R3").Value = 18.50467362471
$ws.Range("S3").Value = 0.004037025134250563
$ws.Range("T3").Value = 0.004037025134250563
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 3.317404
$ws.Range("N4").Value = 9.952211999999999
$ws.Range("O4").Value = 0.1040219928245168
$ws.Range("P4").Value = 0.1040219928245168
$ws.Range("Q4").Value = 1.78623959678
$ws.Range("R4").Value = 16.07615637102
$ws.Range("S4").Value = 0.003507213834092528
$ws.Range("T4").Value = 0.003507213834092528
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("M5").Value = 24.75542533333333
$ws.Range("N5").Value = 74.26627599999999
$ws.Range("O5").Value = 0.7762421087066456
$ws.Range("P5").Value = 0.7762421087066456
$ws.Range("Q5").Value = 298.0659575944183
$ws.Range("R5").Value = 2682.593618349764
$ws.Range("S5").Value = 0.5852412251030922
$ws.Range("T5").Value = 0.5852412251030921
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("M6").Value = 3.818542
$ws.Range("O6").Value = 0.1197358984688377
$ws.Range("P6").Value = 0.1197358984688377
$ws.Range("S6").Value = 0.09027387605328206
$ws.Range("T6").Value = 0.09027387605328205
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 3.317404
$ws.Range("N7").Value = 9.952211999999999
$ws.Range("O7").Value = 0.1040219928245168
$ws.Range("P7").Value = 0.1040219928245168
$ws.Range("Q7").Value = 39.94296953791867
$ws.Range("R7").Value = 359.486725841268
$ws.Range("S7").Value = 0.07842650873413519
$ws.Range("T7").Value = 0.07842650873413519
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("M8").Value = 24.75542533333333
$ws.Range("N8").Value = 74.26627599999999
$ws.Range("O8").Value = 0.7762421087066456
$ws.Range("P8").Value = 0.7762421087066456
$ws.Range("Q8").Value = 83.94816410690888
$ws.Range("R8").Value = 755.5334769621799
$ws.Range("S8").Value = 0.1648290425501541
$ws.Range("T8").Value = 0.164829042550154
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("M9").Value = 3.818542
$ws.Range("O9").Value = 0.1197358984688377
$ws.Range("P9").Value = 0.1197358984688377
$ws.Range("Q9").Value = 12.94906414043667
$ws.Range("R9").Value = 116.54157726393
$ws.Range("S9").Value = 0.02542499728130506
$ws.Range("T9").Value = 0.02542499728130506
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 3.317404
$ws.Range("N10").Value = 9.952211999999999
$ws.Range("O10").Value = 0.1040219928245168
$ws.Range("P10").Value = 0.1040219928245168
$ws.Range("Q10").Value = 11.24965423340666
$ws.Range("R10").Value = 101.24688810066
$ws.Range("S10").Value = 0.02208827025628906
$ws.Range("T10").Value = 0.02208827025628906
